$wb = $excel.ActiveWorkbook

# --- 1. Insert a new worksheet "2022-Q3" right before the existing "2021-Q3" sheet ---
$ws1 = $wb.Worksheets.Item(1)        # "总计" summary sheet
$sheetOld = $wb.Worksheets.Item(2)   # existing "2021-Q3" sheet (will become 3rd sheet)
$new = $wb.Worksheets.Add($sheetOld)
$new.Name = "2022-Q3"

# --- 2. Update the "总计" (summary) sheet: add a new row for 2022-Q3 above the 2021-Q3 row ---
# Preserve the "s=2" direct formatting on column A by copying it onto the row that is
# about to receive the old 2021-Q3 data (row 3).
$ws1.Range("A2").Copy()
$ws1.Range("A3").PasteSpecial(-4122)

# Push the existing 2021-Q3 totals down into row 3
$ws1.Range("A3").Value = 1
$ws1.Range("B3").Value = "2021-Q3"
$ws1.Range("C3").Value = 2
$ws1.Range("D3").Value = 0.08

# Write the new 2022-Q3 totals into row 2
$ws1.Range("A2").Value = 0
$ws1.Range("B2").Value = "2022-Q3"
$ws1.Range("C2").Value = 1
$ws1.Range("D2").Value = 0.11

# --- 3. Populate the new "2022-Q3" sheet with its fund holding data ---
# Copy header/A2 direct formatting ("s=2") from the 总计 sheet so the new sheet matches
# the look of the summary sheet rather than the old 2021-Q3 sheet's formatting.
$ws1.Range("B1:D1").Copy()
$new.Range("B1:D1").PasteSpecial(-4122)
$ws1.Range("B1").Copy()
$new.Range("E1").PasteSpecial(-4122)
$new.Range("F1").PasteSpecial(-4122)
$new.Range("G1").PasteSpecial(-4122)
$new.Range("H1").PasteSpecial(-4122)

$ws1.Range("A2").Copy()
$new.Range("A2").PasteSpecial(-4122)

$new.Range("B1").Value = "基金代码"
$new.Range("C1").Value = "基金名称"
$new.Range("D1").Value = "基金规模"
$new.Range("E1").Value = "股票总仓位"
$new.Range("F1").Value = "仓位占比"
$new.Range("G1").Value = "持有市值(亿元)"
$new.Range("H1").Value = "仓位排名"

$new.Range("A2").Value = 0

# Force these as plain text (matches inlineStr in the source data) so values such as
# the leading-zero fund code and fixed-decimal numbers are not coerced into numbers.
$new.Range("B2:G2").NumberFormat = "@"
$new.Range("B2").Value = "007592"
$new.Range("C2").Value = "华夏价值精选混合"
$new.Range("D2").Value = "2.26"
$new.Range("E2").Value = "93.77"
$new.Range("F2").Value = "4.77"
$new.Range("G2").Value = "0.1078"
$new.Range("B2:G2").ClearFormats()

$new.Range("H2").Value = 8
